$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.710.58"
$ws.Range("E2").Value = "  +1.91%  "

# Row 3
$ws.Range("D3").Value = "1.854.26"
$ws.Range("E3").Value = "  +1.47%  "

# Row 4
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "'244.23"
$ws.Range("E5").Value = "  +0.91%  "

# Row 6
$ws.Range("D6").Value = "'0.6387"
$ws.Range("E6").Value = "  +3.07%  "

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "'46.85"
$ws.Range("E8").Value = "  +3.42%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2984"
$ws.Range("E9").Value = "  +2.66%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07472"
$ws.Range("E10").Value = "  +1.57%  "

# Row 11
$ws.Range("D11").Value = "'24.26"
$ws.Range("E11").Value = "  +5.24%  "

# Row 12
$ws.Range("D12").Value = "'0.07647"
$ws.Range("E12").Value = "  -0.42%  "

# Row 13
$ws.Range("D13").Value = "1.862.78"
$ws.Range("E13").Value = "  +1.84%  "

# Row 14
$ws.Range("D14").Value = "'5.044"
$ws.Range("E14").Value = "  +1.97%  "

# Row 15
$ws.Range("E15").Value = "  +3.45%  "

# Row 16
$ws.Range("D16").Value = "'83.70"
$ws.Range("E16").Value = "  +1.94%  "

# Row 17
$ws.Range("D17").Value = "'0.000009511"
$ws.Range("E17").Value = "  +7.17%  "

# Row 18
$ws.Range("D18").Value = "'6.053"
$ws.Range("E18").Value = "  +3.53%  "

# Row 19
$ws.Range("D19").Value = "29.719.71"
$ws.Range("E19").Value = "  +2.03%  "

# Row 20
$ws.Range("D20").Value = "2.117.86"
$ws.Range("E20").Value = "  +2.29%  "

# Row 21
$ws.Range("D21").Value = "'235.83"
$ws.Range("E21").Value = "  -0.72%  "

# Row 22
$ws.Range("E22").Value = "  +1.46%  "

# Row 23
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
$ws.Range("D24").Value = "'7.398"
$ws.Range("E24").Value = "  +0.98%  "

# Row 25
$ws.Range("E25").Value = "  -0.05%  "

# Row 26
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("E27").Value = "  +0.55%  "

# Row 28
$ws.Range("D28").Value = "'8.490"
$ws.Range("E28").Value = "  +0.07%  "

# Row 29
$ws.Range("D29").Value = "'17.89"
$ws.Range("E29").Value = "  +1.50%  "

# Row 30
$ws.Range("D30").Value = "'0.06239"
$ws.Range("E30").Value = "  +5.68%  "

# Row 31
$ws.Range("D31").Value = "'1.491"
$ws.Range("E31").Value = "  +0.29%  "

# Row 32
$ws.Range("D32").Value = "'1.272"
$ws.Range("E32").Value = "  +5.11%  "

# Row 33
$ws.Range("D33").Value = "'4.146"
$ws.Range("E33").Value = "  +1.75%  "

# Row 34
$ws.Range("D34").Value = "'4.097"
$ws.Range("E34").Value = "  +0.81%  "

# Row 35
$ws.Range("E35").Value = "  +1.27%  "

# Row 36
$ws.Range("D36").Value = "'1.174"
$ws.Range("E36").Value = "  +3.42%  "

# Row 37
$ws.Range("D37").Value = "'0.7271"
$ws.Range("E37").Value = "  -0.23%  "

# Row 38
$ws.Range("D38").Value = "'2.606"
$ws.Range("E38").Value = "  +0.14%  "

# Row 39
$ws.Range("D39").Value = "'2.845"
$ws.Range("E39").Value = "  -0.04%  "

# Row 40
$ws.Range("E40").Value = "  +2.01%  "

# Row 41
$ws.Range("D41").Value = "1.203.39"
$ws.Range("E41").Value = "  -1.16%  "

# Row 42
$ws.Range("D42").Value = "'0.9244"
$ws.Range("E42").Value = "  +0.89%  "

# Row 43
$ws.Range("D43").Value = "'6.143"
$ws.Range("E43").Value = "  -1.97%  "

# Row 44
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("D45").Value = "2.025.13"
$ws.Range("E45").Value = "  +2.57%  "

# Row 46
$ws.Range("D46").Value = "'101.99"
$ws.Range("E46").Value = "  +0.28%  "

# Row 47
$ws.Range("D47").Value = "'66.00"
$ws.Range("E47").Value = "  +2.01%  "

# Row 48
$ws.Range("E48").Value = "  +1.86%  "

# Row 49
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "'0.4055"
$ws.Range("E49").Value = "  +1.05%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.161"
$ws.Range("E50").Value = "  +0.33%  "

# Row 51
$ws.Range("D51").Value = "'0.05793"
$ws.Range("E51").Value = "  +0.80%  "
